$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Remove the "Addable dossier types" / "addable_dossier_types" column (column P)
$ws.Columns.Item(16).Delete()
